$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume cells retain text formatting (values such as "1.001" or
# "25.915.38" must not be auto-converted to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.915.38'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.740.79'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '246.42'
$ws.Range("E5").Value = '  +4.76%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.5060'
$ws.Range("E7").Value = '  -4.27%  '
$ws.Range("D8").Value = '0.2709'
$ws.Range("E8").Value = '  -3.13%  '
$ws.Range("D9").Value = '0.06173'
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").Value = '1.747.14'
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = '0.07230'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '15.06'
$ws.Range("E12").Value = '  -2.58%  '
$ws.Range("D13").Value = '0.6457'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").Value = '4.615'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = '77.44'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '25.939.89'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '11.80'
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").Value = '0.000006794'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '1.966.14'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = '4.276'
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").Value = '8.622'
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("D24").Value = '5.372'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").Value = '136.10'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = '1.501'
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").Value = '15.20'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").Value = '1.768'
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("D29").Value = '105.29'
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D30").Value = '3.899'
$ws.Range("E30").Value = '  +2.55%  '
$ws.Range("D31").Value = '0.08215'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '3.628'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").Value = '0.04670'
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = '2.654'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = '0.9907'
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("D36").Value = '0.6200'
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("D37").Value = '2.724'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").Value = '0.01596'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").Value = '1.909'
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").Value = '1.001'
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").Value = '99.03'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").Value = '0.7573'
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("D43").Value = '0.3842'
$ws.Range("E43").Value = '  -2.26%  '
$ws.Range("D44").Value = '4.974'
$ws.Range("E44").Value = '  -1.03%  '
$ws.Range("D45").Value = '0.1131'
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").Value = '6.249'
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = '55.39'
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").Value = '0.05240'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("D49").Value = '30.63'
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("D50").Value = '7.478'
$ws.Range("E50").Value = '  -1.73%  '
$ws.Range("D51").Value = '0.3398'
$ws.Range("E51").Value = '  -1.62%  '
